# Updated cryptos list values (price + 1h volume change) in the Price/Volume(1h)
# columns of the crypto tracker sheet. Cells that look like plain decimals
# (e.g. "2.38") are forced to Text via NumberFormat "@" before assignment so
# they keep being stored as text (matching the original inlineStr cells)
# instead of Excel auto-coercing them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.884.73'
$ws.Range('E2').Value = '  -1.28%  '
$ws.Range('D3').Value = '2.041.34'
$ws.Range('E3').Value = '  -2.06%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '250.79'
$ws.Range('E5').Value = '  -0.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.665'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '57.62'
$ws.Range('E7').Value = '  +0.41%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '61.30'
$ws.Range('E9').Value = '  -1.89%  '
$ws.Range('E10').Value = '  +0.53%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0782'
$ws.Range('E11').Value = '  +3.70%  '
$ws.Range('E12').Value = '  +1.48%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '16.31'
$ws.Range('E13').Value = '  +5.94%  '
$ws.Range('D14').Value = '2.340.00'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.802'
$ws.Range('E15').Value = '  -5.61%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.53'
$ws.Range('E16').Value = '  +4.74%  '
$ws.Range('D17').Value = '2.049.35'
$ws.Range('E17').Value = '  -1.82%  '
$ws.Range('D18').Value = '36.833.37'
$ws.Range('E18').Value = '  -1.48%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.64'
$ws.Range('E19').Value = '  +12.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '75.23'
$ws.Range('E20').Value = '  +2.68%  '
$ws.Range('E21').Value = '  +6.23%  '
$ws.Range('E22').Value = '  +2.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.87'
$ws.Range('E23').Value = '  -1.69%  '
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.38'
$ws.Range('E25').Value = '  -3.92%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.34'
$ws.Range('E26').Value = '  +15.88%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '168.89'
$ws.Range('E27').Value = '  -1.64%  '
$ws.Range('E28').Value = '  -0.71%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.12'
$ws.Range('E29').Value = '  -4.69%  '
$ws.Range('E30').Value = '  +0.59%  '
$ws.Range('E31').Value = '  +5.22%  '
$ws.Range('E32').Value = '  +2.69%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0618'
$ws.Range('E33').Value = '  -1.56%  '
$ws.Range('E34').Value = '  +2.43%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0882'
$ws.Range('E35').Value = '  -3.26%  '
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('E37').Value = '  -3.71%  '
$ws.Range('E38').Value = '  -4.69%  '
$ws.Range('E39').Value = '  +12.61%  '
$ws.Range('E40').Value = '  -0.42%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '17.69'
$ws.Range('E41').Value = '  -0.61%  '
$ws.Range('E42').Value = '  -2.47%  '
$ws.Range('E43').Value = '  -3.69%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '96.50'
$ws.Range('E44').Value = '  -3.74%  '
$ws.Range('E45').Value = '  +1.09%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.68'
$ws.Range('E46').Value = '  +15.61%  '
$ws.Range('E47').Value = '  +4.54%  '
$ws.Range('D48').Value = '1.278.91'
$ws.Range('E48').Value = '  -3.73%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.89'
$ws.Range('E49').Value = '  -2.08%  '
$ws.Range('E50').Value = '  -4.69%  '
$ws.Range('D51').Value = '2.232.75'
$ws.Range('E51').Value = '  -1.90%  '
